# fig_archaeology / arch_strat_ROAD consolidation update
# Renames 6 top-level archaeological-period labels in column A (and, where
# column B mirrors column A, in column B too) to their "Other ..." variants,
# and for the 3 "Stone Age" categories (ESA/LSA/MSA) also clears the cell's
# fill formatting back to the worksheet default.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arch_strat_ROAD")

# Rows whose column-A cell currently reads "ESA"/"LSA"/"MSA" -> renamed to
# "Other Early/Late/Middle Stone Age" AND lose their highlight fill.
$rowsClearFill = @{
    "ESA" = @(15, 19, 20, 21, 29);
    "LSA" = @(56, 57, 58, 59, 60, 61, 62, 63, 64, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77);
    "MSA" = @(104, 105, 106, 107, 109, 111, 113, 114, 115, 118, 119, 120, 122, 123, 124);
}
$newTextClearFill = @{
    "ESA" = "Other Early Stone Age";
    "LSA" = "Other Late Stone Age";
    "MSA" = "Other Middle Stone Age";
}

# Rows whose column-A cell currently reads "Lower/Middle/Upper Paleolithic"
# -> renamed to "Other Lower/Middle/Upper Paleolithic"; their highlight
# fill is left exactly as-is.
$rowsKeepFill = @{
    "Lower Paleolithic"  = @(34, 36, 38, 42, 43, 44, 45, 46, 48, 51, 53);
    "Middle Paleolithic" = @(79, 80, 82, 84, 85, 86, 88, 93, 94, 95, 96, 97);
    "Upper Paleolithic"  = @(136, 137, 139, 140, 141, 144, 146, 147, 148, 151, 152, 153, 154, 155);
}
$newTextKeepFill = @{
    "Lower Paleolithic"  = "Other Lower Paleolithic";
    "Middle Paleolithic" = "Other Middle Paleolithic";
    "Upper Paleolithic"  = "Other Upper Paleolithic";
}

# Old label -> new label, used to also re-sync column B whenever B mirrors A.
$oldToNew = @{
    "ESA" = "Other Early Stone Age";
    "LSA" = "Other Late Stone Age";
    "MSA" = "Other Middle Stone Age";
    "Lower Paleolithic"  = "Other Lower Paleolithic";
    "Middle Paleolithic" = "Other Middle Paleolithic";
    "Upper Paleolithic"  = "Other Upper Paleolithic";
}

foreach ($label in $rowsClearFill.Keys) {
    $newText = $newTextClearFill[$label]
    foreach ($r in $rowsClearFill[$label]) {
        $ws.Cells.Item($r, 1).Value2 = $newText
        $ws.Cells.Item($r, 1).ClearFormats()
    }
}

foreach ($label in $rowsKeepFill.Keys) {
    $newText = $newTextKeepFill[$label]
    foreach ($r in $rowsKeepFill[$label]) {
        $ws.Cells.Item($r, 1).Value2 = $newText
    }
}

# Column B: wherever it currently mirrors one of the six renamed labels,
# update it to the corresponding new label too.
$lastRow = 155
for ($r = 1; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($null -ne $bVal -and $oldToNew.ContainsKey($bVal)) {
        $bCell.Value2 = $oldToNew[$bVal]
    }
}

# Update the saved view state to match (scrolled/selected position).
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("B155").Select()
